$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1.02
$ws.Range("C2").Value = 1.037267621666352
$ws.Range("D2").Value = 1.041567248677424
$ws.Range("E2").Value = 1.040917735927591
$ws.Range("F2").Value = 1.047099851320691
$ws.Range("I2").Value = 1.041759102693556
$ws.Range("J2").Value = 1.042371420824399
$ws.Range("K2").Value = 1.04434639827114
$ws.Range("L2").Value = 1.043698724709446
$ws.Range("M2").Value = 1.049863434889588
$ws.Range("N2").Value = 1.043851707803643

$ws.Range("B3").Value = 1.02
$ws.Range("C3").Value = 1.038273741405853
$ws.Range("D3").Value = 1.042356710996973
$ws.Range("E3").Value = 1.041872374193163
$ws.Range("F3").Value = 1.04835217711232
$ws.Range("I3").Value = 1.04208999153714
$ws.Range("J3").Value = 1.043021564316282
$ws.Range("K3").Value = 1.044946570257686
$ws.Range("L3").Value = 1.044463505605055
$ws.Range("M3").Value = 1.050926393851624
$ws.Range("N3").Value = 1.044502774573857

$ws.Range("B4").Value = 1.02
$ws.Range("C4").Value = 1.038924809497708
$ws.Range("D4").Value = 1.042867503925206
$ws.Range("E4").Value = 1.042490473463431
$ws.Range("F4").Value = 1.049162899242599
$ws.Range("I4").Value = 1.042302795913146
$ws.Range("J4").Value = 1.043441700525213
$ws.Range("K4").Value = 1.04533422457347
$ws.Range("L4").Value = 1.044958134474317
$ws.Range("M4").Value = 1.051614025815972
$ws.Range("N4").Value = 1.044923507424394

$ws.Range("B5").Value = 1.02
$ws.Range("C5").Value = 1.039198528522738
$ws.Range("D5").Value = 1.043082230735498
$ws.Range("E5").Value = 1.042750414120897
$ws.Range("F5").Value = 1.049503819679426
$ws.Range("I5").Value = 1.042391946892639
$ws.Range("J5").Value = 1.043618193947927
$ws.Range("K5").Value = 1.045497027336548
$ws.Range("L5").Value = 1.04516601985664
$ws.Range("M5").Value = 1.051903065519271
$ws.Range("N5").Value = 1.045100251488012

$ws.Range("B6").Value = 1.02
$ws.Range("C6").Value = 1.039244487740846
$ws.Range("D6").Value = 1.043118283679765
$ws.Range("E6").Value = 1.04279406467054
$ws.Range("F6").Value = 1.049561067195558
$ws.Range("I6").Value = 1.042406897437425
$ws.Range("J6").Value = 1.043647820234926
$ws.Range("K6").Value = 1.045524352830787
$ws.Range("L6").Value = 1.045200921402188
$ws.Range("M6").Value = 1.051951594188922
$ws.Range("N6").Value = 1.045129919847738

$ws.Range("B7").Value = 1.02
$ws.Range("C7").Value = 1.038928466905292
$ws.Range("D7").Value = 1.042870373157369
$ws.Range("E7").Value = 1.042493946443204
$ws.Range("F7").Value = 1.049167454272027
$ws.Range("I7").Value = 1.042303988378863
$ws.Range("J7").Value = 1.043444059357183
$ws.Range("K7").Value = 1.045336400608875
$ws.Range("L7").Value = 1.044960912471613
$ws.Range("M7").Value = 1.051617888137871
$ws.Range("N7").Value = 1.044925869606177

$ws.Range("B8").Value = 1.02
$ws.Range("C8").Value = 1.037607636043433
$ws.Range("D8").Value = 1.041834059565744
$ws.Range("E8").Value = 1.04124028059802
$ws.Range("F8").Value = 1.047523001869569
$ws.Range("I8").Value = 1.041871197839424
$ws.Range("J8").Value = 1.04259125367808
$ws.Range("K8").Value = 1.044549373405865
$ws.Range("L8").Value = 1.043957234629909
$ws.Range("M8").Value = 1.050222703251344
$ws.Range("N8").Value = 1.044071852845191

$ws.Range("B9").Value = 1.02
$ws.Range("C9").Value = 1.035280474948871
$ws.Range("D9").Value = 1.040007642124744
$ws.Range("E9").Value = 1.039034119551986
$ws.Range("F9").Value = 1.044628158790846
$ws.Range("I9").Value = 1.041098590755479
$ws.Range("J9").Value = 1.04108429392576
$ws.Range("K9").Value = 1.043157201643976
$ws.Range("L9").Value = 1.042186831305602
$ws.Range("M9").Value = 1.047762844933195
$ws.Range("N9").Value = 1.042562753037172

$ws.Range("B10").Value = 1.02
$ws.Range("C10").Value = 1.033729233830562
$ws.Range("D10").Value = 1.038789854747535
$ws.Range("E10").Value = 1.037565350505097
$ws.Range("F10").Value = 1.042700141355292
$ws.Range("I10").Value = 1.040576815809375
$ws.Range("J10").Value = 1.040076831311573
$ws.Range("K10").Value = 1.042225515024016
$ws.Range("L10").Value = 1.041005365573759
$ws.Range("M10").Value = 1.04612197068496
$ws.Range("N10").Value = 1.041553859710515

$ws.Range("B11").Value = 1.02
$ws.Range("C11").Value = 1.033057573511303
$ws.Range("D11").Value = 1.03826250254461
$ws.Range("E11").Value = 1.036929835392242
$ws.Range("F11").Value = 1.041865718314382
$ws.Range("I11").Value = 1.040349291642533
$ws.Range("J11").Value = 1.03963991989219
$ws.Range("K11").Value = 1.041821238841577
$ws.Range("L11").Value = 1.040493495312337
$ws.Range("M11").Value = 1.045411212021718
$ws.Range("N11").Value = 1.041116327826802

$ws.Range("B12").Value = 1.02
$ws.Range("C12").Value = 1.032808094190667
$ws.Range("D12").Value = 1.038066614421357
$ws.Range("E12").Value = 1.036693847741903
$ws.Range("F12").Value = 1.041555838281635
$ws.Range("I12").Value = 1.040264539768276
$ws.Range("J12").Value = 1.039477530310691
$ws.Range("K12").Value = 1.041670945053638
$ws.Range("L12").Value = 1.040303320680062
$ws.Range("M12").Value = 1.045147165999383
$ws.Range("N12").Value = 1.040953707633471

$ws.Range("B13").Value = 1.02
$ws.Range("C13").Value = 1.032861608121339
$ws.Range("D13").Value = 1.038108633356879
$ws.Range("E13").Value = 1.036744464681185
$ws.Range("F13").Value = 1.041622305825531
$ws.Range("I13").Value = 1.04028273016138
$ws.Range("J13").Value = 1.039512368008654
$ws.Range("K13").Value = 1.041703189355334
$ws.Range("L13").Value = 1.040344115740117
$ws.Range("M13").Value = 1.045203806533145
$ws.Range("N13").Value = 1.040988594804961

$ws.Range("B14").Value = 1.02
$ws.Range("C14").Value = 1.033036951368379
$ws.Range("D14").Value = 1.038246310492553
$ws.Range("E14").Value = 1.036910327131275
$ws.Range("F14").Value = 1.041840102281109
$ws.Range("I14").Value = 1.040342290907522
$ws.Range("J14").Value = 1.039626498783995
$ws.Range("K14").Value = 1.041808818113175
$ws.Range("L14").Value = 1.040477776299633
$ws.Range("M14").Value = 1.045389386689416
$ws.Range("N14").Value = 1.041102887659094

$ws.Range("B15").Value = 1.02
$ws.Range("C15").Value = 1.03314498681133
$ws.Range("D15").Value = 1.038331137129349
$ws.Range("E15").Value = 1.037012529866946
$ws.Range("F15").Value = 1.041974302016362
$ws.Range("I15").Value = 1.040378956541775
$ws.Range("J15").Value = 1.039696805106524
$ws.Range("K15").Value = 1.041873882568509
$ws.Range("L15").Value = 1.040560123254209
$ws.Range("M15").Value = 1.045503723598107
$ws.Range("N15").Value = 1.041173293824665

$ws.Range("B16").Value = 1.02
$ws.Range("C16").Value = 1.033773810467776
$ws.Range("D16").Value = 1.038824852513186
$ws.Range("E16").Value = 1.037607537495085
$ws.Range("F16").Value = 1.04275552801563
$ws.Range("I16").Value = 1.040591882288414
$ws.Range("J16").Value = 1.040105813490498
$ws.Range("K16").Value = 1.042252327604162
$ws.Range("L16").Value = 1.041039330654637
$ws.Range("M16").Value = 1.046169136045732
$ws.Range("N16").Value = 1.041582883047459

$ws.Range("B17").Value = 1.02
$ws.Range("C17").Value = 1.034168264658005
$ws.Range("D17").Value = 1.039134536042661
$ws.Range("E17").Value = 1.03798089639273
$ws.Range("F17").Value = 1.043245682198515
$ws.Range("I17").Value = 1.040725018669984
$ws.Range("J17").Value = 1.040362193306577
$ws.Range("K17").Value = 1.04248948871499
$ws.Range("L17").Value = 1.041339847752808
$ws.Range("M17").Value = 1.046586464137044
$ws.Range("N17").Value = 1.041839626952284

$ws.Range("B18").Value = 1.02
$ws.Range("C18").Value = 1.034398346739913
$ws.Range("D18").Value = 1.039315165110485
$ws.Range("E18").Value = 1.038198715814017
$ws.Range("F18").Value = 1.043531621580401
$ws.Range("I18").Value = 1.04080252125148
$ws.Range("J18").Value = 1.040511670339823
$ws.Range("K18").Value = 1.042627738788852
$ws.Range("L18").Value = 1.04151510652216
$ws.Range("M18").Value = 1.046829860577049
$ws.Range("N18").Value = 1.041989316260063

$ws.Range("B19").Value = 1.02
$ws.Range("C19").Value = 1.034476799410652
$ws.Range("D19").Value = 1.039376754250602
$ws.Range("E19").Value = 1.038272994320768
$ws.Range("F19").Value = 1.043629126474535
$ws.Range("I19").Value = 1.040828921574663
$ws.Range("J19").Value = 1.040562627125729
$ws.Range("K19").Value = 1.042674864582985
$ws.Range("L19").Value = 1.041574860522798
$ws.Range("M19").Value = 1.046912848441777
$ws.Range("N19").Value = 1.042040345410449

$ws.Range("B20").Value = 1.02
$ws.Range("C20").Value = 1.034125943083983
$ws.Range("D20").Value = 1.039101310352048
$ws.Range("E20").Value = 1.037940833809989
$ws.Range("F20").Value = 1.043193089089798
$ws.Range("I20").Value = 1.040710750285986
$ws.Range("J20").Value = 1.040334692899674
$ws.Range("K20").Value = 1.042464052052405
$ws.Range("L20").Value = 1.041307607984394
$ws.Range("M20").Value = 1.046541691270406
$ws.Range("N20").Value = 1.041812087491648

$ws.Range("B21").Value = 1.02
$ws.Range("C21").Value = 1.032985317011172
$ws.Range("D21").Value = 1.038205768171716
$ws.Range("E21").Value = 1.036861482814045
$ws.Range("F21").Value = 1.041775964957979
$ws.Range("I21").Value = 1.040324758357764
$ws.Range("J21").Value = 1.039592892910027
$ws.Range("K21").Value = 1.041777716599084
$ws.Range("L21").Value = 1.040438417789407
$ws.Range("M21").Value = 1.045334739033908
$ws.Range("N21").Value = 1.041069234060931

$ws.Range("B22").Value = 1.02
$ws.Range("C22").Value = 1.032268190186737
$ws.Range("D22").Value = 1.037642670496048
$ws.Range("E22").Value = 1.036183262603815
$ws.Range("F22").Value = 1.040885319827425
$ws.Range("I22").Value = 1.04008068555062
$ws.Range("J22").Value = 1.039125908123932
$ws.Range("K22").Value = 1.041345451951164
$ws.Range("L22").Value = 1.039891673328912
$ws.Range("M22").Value = 1.044575656312803
$ws.Range("N22").Value = 1.040601586102875

$ws.Range("B23").Value = 1.02
$ws.Range("C23").Value = 1.032648349961956
$ws.Range("D23").Value = 1.037941182420188
$ws.Range("E23").Value = 1.036542761012983
$ws.Range("F23").Value = 1.041357434300911
$ws.Range("I23").Value = 1.040210204405767
$ws.Range("J23").Value = 1.039373521127186
$ws.Range("K23").Value = 1.041574673586707
$ws.Range("L23").Value = 1.040181536613262
$ws.Range("M23").Value = 1.04497808203865
$ws.Range("N23").Value = 1.040849550744996

$ws.Range("B24").Value = 1.02
$ws.Range("C24").Value = 1.034145066370162
$ws.Range("D24").Value = 1.039116323624043
$ws.Range("E24").Value = 1.03795893622598
$ws.Range("F24").Value = 1.04321685352547
$ws.Range("I24").Value = 1.040717198029103
$ws.Range("J24").Value = 1.040347119350459
$ws.Range("K24").Value = 1.042475546038556
$ws.Range("L24").Value = 1.041322175833683
$ws.Range("M24").Value = 1.046561922275323
$ws.Range("N24").Value = 1.041824531589419

$ws.Range("B25").Value = 1.02
$ws.Range("C25").Value = 1.035882065967686
$ws.Range("D25").Value = 1.04047984770424
$ws.Range("E25").Value = 1.039604112835807
$ws.Range("F25").Value = 1.045376209244144
$ws.Range("I25").Value = 1.041299509973472
$ws.Range("J25").Value = 1.041474377107071
$ws.Range("K25").Value = 1.043517741824863
$ws.Range("L25").Value = 1.042644734670666
$ws.Range("M25").Value = 1.048398944249698
$ws.Range("N25").Value = 1.042953390181344

